$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.715.23"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").Value = "2.312.15"
$ws.Range("E3").Value = "  +0.73%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.75"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.51"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.84%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.502"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.47%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -1.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "34.17"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.57%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "18.91"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.47%  "
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("E13").Value = "  +0.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.73"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.67%  "
$ws.Range("D15").Value = "2.671.75"
$ws.Range("E15").Value = "  +0.69%  "
$ws.Range("D16").Value = "2.307.96"
$ws.Range("E16").Value = "  +0.71%  "
$ws.Range("E17").Value = "  +1.34%  "
$ws.Range("D18").Value = "42.666.80"
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.13"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -4.73%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.13"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.95%  "
$ws.Range("D21").Value = "0.0₃0890"
$ws.Range("E21").Value = "  -0.25%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.71"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.93%  "
$ws.Range("E23").Value = "  +5.57%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "235.08"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.21%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.42"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.73%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.30"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.22%  "
$ws.Range("E28").Value = "  +14.48%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "165.78"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.29%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.10"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "32.07"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.34%  "
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.00"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.76%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.70"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.46"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.69%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0698"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.76%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.33"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.75%  "
$ws.Range("E38").Value = "  +3.02%  "
$ws.Range("E39").Value = "  -0.29%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.72"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.79%  "
$ws.Range("E41").Value = "  -0.70%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.49"
$ws.Range("D42").ClearFormats()
$ws.Range("D43").Value = "1.923.66"
$ws.Range("E43").Value = "  -3.76%  "
$ws.Range("E44").Value = "  -0.58%  "
$ws.Range("E45").Value = "  -1.99%  "
$ws.Range("E46").Value = "  -1.49%  "
$ws.Range("E47").Value = "  -0.40%  "
$ws.Range("E48").Value = "  +2.24%  "
$ws.Range("D49").Value = "2.540.71"
$ws.Range("E49").Value = "  +0.79%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.31"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.24%  "
$ws.Range("E51").Value = "  +1.54%  "
